$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# Row 3: fill in CRM value, Batch value %, D3 and E3, F3 uses existing "No open date" note
$ws.Range("A3").Value = 20241009
$ws.Range("B3").Value = 2221.9117900000001
$ws.Range("C3").Value = 2224.4699999999998
$ws.Range("D3").Formula = "=100*(B3-C3)/C3"
$ws.Range("E3").Value = 180
$ws.Range("F3").Value = "No open date"

# Row 4
$ws.Range("A4").Value = 20241010
$ws.Range("B4").Value = 2319.57584658833
$ws.Range("C4").Value = 2224.4699999999998
$ws.Range("D4").Formula = "=100*(B4-C4)/C4"
$ws.Range("E4").Value = 180
$ws.Range("F4").Value = "open 20241010"

# Row 5
$ws.Range("A5").Value = 20241031
$ws.Range("B5").Value = 2199.6047899999999
$ws.Range("C5").Value = 2215.13
$ws.Range("D5").Formula = "=100*(B5-C5)/C5"
$ws.Range("E5").Value = 202
$ws.Range("F5").Value = "CRM202_opened20241031"

$ws.Range("F9").Select()
